$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Unique Concord, Tata Nagar, Bangalore" address (row 3) was removed
# from the list; deleting the whole row shifts the rows below it up by
# one, which also drops the now-unused shared string and shrinks the
# sheet's used range from A1:A5 down to A1:A4.
$ws.Range("A3").EntireRow.Delete()

# Leave the selection on the new last row, matching the saved file.
$ws.Range("A4").Select()
